$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 5 new "shirts" sale rows (rows 8-12), matching the existing
# sale pattern used by update_stock_with_sale / update_sales for this
# product: year, month, day, code product, name, amount, price.
for ($i = 8; $i -le 12; $i++) {
    $ws.Cells.Item($i, 1).Value = 2018
    $ws.Cells.Item($i, 2).Value = 12
    $ws.Cells.Item($i, 3).Value = 27
    $ws.Cells.Item($i, 4).Value = 1
    $ws.Cells.Item($i, 5).Value = "shirts"
    $ws.Cells.Item($i, 6).Value = 2
    $ws.Cells.Item($i, 7).Value = 89.90000000000001
}
